$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.988.74'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").Value = '3.333.62'
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '584.57'
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("D6").Value = '176.72'
$ws.Range("E6").Value = '  +2.07%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("E8").Value = '  +1.46%  '
$ws.Range("E9").Value = '  +4.98%  '
$ws.Range("D10").Value = '0.583'
$ws.Range("E10").Value = '  +1.66%  '
$ws.Range("D11").Value = '47.86'
$ws.Range("E11").Value = '  +5.85%  '
$ws.Range("E12").Value = '  +2.27%  '
$ws.Range("D13").Value = '693.29'
$ws.Range("E13").Value = '  +4.97%  '
$ws.Range("D14").Value = '3.877.68'
$ws.Range("E14").Value = '  +0.59%  '
$ws.Range("D15").Value = '8.43'
$ws.Range("E15").Value = '  +1.03%  '
$ws.Range("D16").Value = '68.033.81'
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D18").Value = '3.314.03'
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("D19").Value = '17.48'
$ws.Range("E19").Value = '  +0.62%  '
$ws.Range("D20").Value = '11.14'
$ws.Range("E20").Value = '  +2.78%  '
$ws.Range("D21").Value = '0.894'
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("D23").Value = '16.92'
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("D24").Value = '100.51'
$ws.Range("E24").Value = '  +3.79%  '
$ws.Range("D25").Value = '3.91'
$ws.Range("E25").Value = '  +2.06%  '
$ws.Range("D26").Value = '2.69'
$ws.Range("E26").Value = '  +1.23%  '
$ws.Range("D27").Value = '9.47'
$ws.Range("E27").Value = '  +2.29%  '
$ws.Range("D28").Value = '33.03'
$ws.Range("E28").Value = '  -0.58%  '
$ws.Range("D29").Value = '8.52'
$ws.Range("E29").Value = '  +2.02%  '
$ws.Range("D30").Value = '6.96'
$ws.Range("E30").Value = '  -3.97%  '
$ws.Range("D31").Value = '570.58'
$ws.Range("E31").Value = '  -3.28%  '
$ws.Range("D32").Value = '11.04'
$ws.Range("E32").Value = '  +1.18%  '
$ws.Range("E33").Value = '  +2.02%  '
$ws.Range("D34").Value = '3.740.80'
$ws.Range("E34").Value = '  +0.83%  '
$ws.Range("D35").Value = '57.49'
$ws.Range("E35").Value = '  +1.11%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").Value = '3.32'
$ws.Range("E37").Value = '  +1.42%  '
$ws.Range("D38").Value = '35.24'
$ws.Range("E38").Value = '  +9.07%  '
$ws.Range("E39").Value = '  +3.07%  '
$ws.Range("D40").Value = '3.16'
$ws.Range("E40").Value = '  +2.78%  '
$ws.Range("D41").Value = '2.62'
$ws.Range("E41").Value = '  +0.62%  '
$ws.Range("D42").Value = '0.0₃0674'
$ws.Range("E42").Value = '  +2.31%  '
$ws.Range("D43").Value = '0.334'
$ws.Range("E43").Value = '  +1.12%  '
$ws.Range("E44").Value = '  -1.14%  '
$ws.Range("E45").Value = '  +1.08%  '
$ws.Range("D46").Value = '2.64'
$ws.Range("E46").Value = '  +2.39%  '
$ws.Range("E47").Value = '  +1.35%  '
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("E49").Value = '  -0.82%  '
$ws.Range("D50").Value = '130.96'
$ws.Range("E50").Value = '  +2.70%  '
$ws.Range("D51").Value = '2.57'
$ws.Range("E51").Value = '  -1.10%  '
